# edit.ps1
# Applies the "Saldo" worksheet update described by the commit diff:
#  - Row 2 (account 005645211/AGUINALDO/150000) is replaced in place with
#    005053939/VIRGILIO/193497.83
#  - Several new account rows are inserted at various points in the table
#  - A few trailing (negative-balance) rows are removed
#
# NOTE: "Conta" (account number) values are zero-padded numeric-looking
# strings (e.g. "005053939") that must be preserved as TEXT, not coerced
# into numbers (which would drop the leading zeros). We force text
# formatting on those cells before assigning the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-AccountCell($row, $value) {
    # Column A holds account numbers as text (e.g. "005053939").
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Set-DataRow($row, $conta, $nome, $saldo) {
    Set-AccountCell $row $conta
    $ws.Cells.Item($row, 2).Value = $nome
    $ws.Cells.Item($row, 3).Value = $saldo
}

# ---------------------------------------------------------------------
# 1) Replace the first data row (Excel row 2) in place.
# ---------------------------------------------------------------------
Set-DataRow 2 "005053939" "VIRGILIO" 193497.83

# ---------------------------------------------------------------------
# 2) Insert five new rows right before the row that holds account
#    004392159 (RODRIGO / 900.21), originally row 5.
# ---------------------------------------------------------------------
$ws.Range("5:9").Insert()
Set-DataRow 5 "004690692" "PHYLIA"    15273.35
Set-DataRow 6 "004479734" "RODRIGO"   6976.99
Set-DataRow 7 "000772433" "MARCELO"   4977.32
Set-DataRow 8 "004267044" "PATRICIA"  2490.72
Set-DataRow 9 "004862746" "CESAR"     1436.09

# ---------------------------------------------------------------------
# 3) Insert one new row right before the row that holds account
#    005547702 (NATHALIA / 99.9). After step 2 this row is now at 12.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Insert()
Set-DataRow 12 "001761119" "BLUEMETRIX" 213.42

# ---------------------------------------------------------------------
# 4) Insert one new row right before the row that holds account
#    005092207 (BRUNO / 11.63). After steps 2-3 this row is now at 142.
# ---------------------------------------------------------------------
$ws.Rows.Item(142).Insert()
Set-DataRow 142 "005324840" "PEDRO" 12.01

# ---------------------------------------------------------------------
# 5) Remove the old trailing row for account 004690692 (-1689.68), which
#    after the previous insertions now sits at row 232.
# ---------------------------------------------------------------------
$ws.Rows.Item(232).Delete()

# ---------------------------------------------------------------------
# 6) Remove the old trailing rows for accounts 005324840 (-5994.77) and
#    004862746 (-9072.28). These are now adjacent at row 234; deleting
#    row 234 twice removes both.
# ---------------------------------------------------------------------
$ws.Rows.Item(234).Delete()
$ws.Rows.Item(234).Delete()
